# Gespreksformulier Associate Degree.docx - content edit
#
# Summary of changes (per the commit diff):
#  1. After the empty paragraph following "... adviseren van klanten op de
#     shop etc." (just before the "Beschrijf de opdracht" heading), insert
#     several new paragraphs describing a "Platform als service" / "Closed
#     source" idea.
#  2. After the empty paragraph following "... architectuurontwerp en
#     databasebeheer ..." (just before the "Welke begeleiding" heading),
#     insert new paragraphs about "Vue js" / "Remix shopify" tooling.
#  3. Merge three runs (incl. a stray lastRenderedPageBreak) in the
#     "Ik kan wekelijks rekenen ..." paragraph into a single run.
#  4. Split the run in table cell [5,2] ("Ik ga actief vragen ...") into two
#     runs with a lastRenderedPageBreak in between (repagination artifact
#     caused by the new content added above).
#  5. Add a lastRenderedPageBreak before the "6" in table cell [6,1] (same
#     repagination artifact).

$d = $word.ActiveDocument

$pkgOpen = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>"
$pkgClose = "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

# ---------------------------------------------------------------------
# Edit 1: insert the "Platform als service / Closed source" block
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("adviseren van klanten op de shop etc.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorPara = $rng.Paragraphs(1)
$emptyPara = $anchorPara.Next()
$s = $emptyPara.Range.Start
$e = $emptyPara.Range.End
$target = $d.Range($s, $e)

$body1 = @"
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:r>
    <w:t>Platform als service.</w:t>
  </w:r>
  <w:r>
    <w:br/>
    <w:t>Als je saas gewend bent kun je inhaken op symfony.</w:t>
  </w:r>
  <w:r>
    <w:br/>
    <w:t>bevoorbeeld als een ai die prijzen aanpast op basis van de ingelogde klant (tenis kleden is met korting voor tennis leden)</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:t>Closed source.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:t>Laravael/symfony api die extern gehost wordt en gebruikt wordt door een saas programma. Die haalt de data</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorBidi"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:t>Gesloten platform.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
    </w:rPr>
  </w:pPr>
</w:p>
"@

$target.InsertXML($pkgOpen + $body1 + $pkgClose)

# ---------------------------------------------------------------------
# Edit 2: insert the "Vue js / Remix shopify" block
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("architectuurontwerp en databasebeheer", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorPara = $rng.Paragraphs(1)
$emptyPara = $anchorPara.Next()
$s = $emptyPara.Range.Start
$e = $emptyPara.Range.End
$target = $d.Range($s, $e)

$body2 = @"
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
      <w:szCs w:val="22"/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
      <w:spacing w:val="-2"/>
      <w:szCs w:val="22"/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
    <w:t xml:space="preserve">Vue </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
      <w:spacing w:val="-2"/>
      <w:szCs w:val="22"/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
    <w:t>js</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
      <w:szCs w:val="22"/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
      <w:szCs w:val="22"/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
    <w:t xml:space="preserve">Remix </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
      <w:szCs w:val="22"/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
    <w:t>shopify</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
      <w:szCs w:val="22"/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
    <w:t xml:space="preserve"> en gebruik t</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
      <w:szCs w:val="22"/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
    <w:t>ypescript. In samen werking.</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
      <w:szCs w:val="22"/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
    <w:br/>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
      <w:szCs w:val="22"/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
    <w:t>symfony</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
      <w:szCs w:val="22"/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
    <w:t xml:space="preserve"> met </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
      <w:szCs w:val="22"/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
    <w:t>graphql</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
      <w:szCs w:val="22"/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
  </w:pPr>
</w:p>
"@

$target.InsertXML($pkgOpen + $body2 + $pkgClose)

# ---------------------------------------------------------------------
# Edit 3: merge the three runs of the "Ik kan wekelijks rekenen ..."
# paragraph (and drop the lastRenderedPageBreak that used to fall there)
# into a single run with the full sentence.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("overleggen met ervaren professionals, waaronder senior developers. Deze meetings bieden niet alleen de gelegenheid om de voortgang van het afstudeerproject te bespreken.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s = $rng.Start
$e = $rng.End
$target = $d.Range($s, $e)
$target.Text = "overleggen met ervaren professionals, waaronder senior developers. Deze meetings bieden niet alleen de gelegenheid om de voortgang van het afstudeerproject te bespreken."

# ---------------------------------------------------------------------
# Edit 4: split "Ik ga actief vragen ..." (table cell, leeruitkomst 5)
# into two runs with a lastRenderedPageBreak in between.
# ---------------------------------------------------------------------
$t = $d.Tables(1)
$cell = $t.Cell(5, 2)
$pr = $cell.Range.Paragraphs(1).Range
$s = $pr.Start
$e = $pr.End
$target = $d.Range($s, $e)

$body4 = @"
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
      <w:sz w:val="18"/>
      <w:szCs w:val="18"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
    </w:rPr>
    <w:t xml:space="preserve">Ik ga actief vragen om feedback van mijn teamleden en begeleiders, en ik zal die feedback gebruiken om beter te worden in wat ik doe. Door open te staan voor </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>suggesties en positieve kritiek, denk ik dat ik steeds beter zal worden in leren en mijn aanpak steeds verder zal verbeteren.</w:t>
  </w:r>
</w:p>
"@

$target.InsertXML($pkgOpen + $body4 + $pkgClose)

# ---------------------------------------------------------------------
# Edit 5: add a lastRenderedPageBreak before the "6" (table cell,
# leeruitkomst 6).
# ---------------------------------------------------------------------
$t = $d.Tables(1)
$cell = $t.Cell(6, 1)
$pr = $cell.Range.Paragraphs(1).Range
$s = $pr.Start
$e = $pr.End
$target = $d.Range($s, $e)

$body5 = @"
<w:p>
  <w:pPr>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
      <w:sz w:val="18"/>
      <w:szCs w:val="18"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
      <w:sz w:val="18"/>
      <w:szCs w:val="18"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>6</w:t>
  </w:r>
</w:p>
"@

$target.InsertXML($pkgOpen + $body5 + $pkgClose)

Write-Host "All edits applied."
